$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.731.69'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '1.776.69'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '223.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.19'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.289'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.85%  '
$ws.Range('E10').Value = '  -5.34%  '
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('D12').Value = '2.031.84'
$ws.Range('E12').Value = '  -0.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.28%  '
$ws.Range('D14').Value = '1.791.97'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').Value = '33.748.83'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.609'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.13'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.51'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '238.26'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.99%  '
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.59'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.61%  '
$ws.Range('E23').Value = '  -1.78%  '
$ws.Range('E24').Value = '  -1.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.06'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.95%  '
$ws.Range('E26').Value = '  -2.18%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0511'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.59'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.93%  '
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.80'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.03%  '
$ws.Range('D35').Value = '1.384.89'
$ws.Range('E35').Value = '  -1.72%  '
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('E37').Value = '  -2.33%  '
$ws.Range('E38').Value = '  -1.14%  '
$ws.Range('E39').Value = '  +5.70%  '
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.909'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '78.12'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.33%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.58'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +14.82%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.66'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.62%  '
$ws.Range('E45').Value = '  +3.39%  '
$ws.Range('E46').Value = '  +1.13%  '
$ws.Range('E47').Value = '  +11.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '107.43'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.70%  '
$ws.Range('E49').Value = '  -1.76%  '
$ws.Range('D50').Value = '1.931.38'
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.10%  '
